# Adapt column header formatting to respective input file names:
#  - "<Header>_old" -> "<Header>_FV2210"
#  - "<Header>_new" -> "<Header>_FV2304"
# and turn the worksheet's used range into a native Excel Table, plus
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row cells (A1:J1 = "_old" suffix, K1 = "diff",
#        L1:U1 = "_new" suffix) ------------------------------------------------
$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($oldHeaders[$i])_FV2210"
}

# Column K (11) stays "diff"
$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($oldHeaders[$i])_FV2304"
}

# --- 2. Turn A1:U64 into a native Excel Table ("Table1") ----------------------
$tableRange = $ws.Range("A1:U64")
$listObject = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# --- 3. Freeze the header row --------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
